$d = $word.ActiveDocument

# Fix the typo "algoRythms" -> "algoRhythms" in the title paragraph.
$d.Content.Find.Execute("algoRythms", $true, $false, $false, $false, $false, $true, 1, $false, "algoRhythms", 2)

# The "_GoBack" bookmark used to sit at the end of the "3.    React" line;
# relocate it so it sits right between "algo" and "Rhythms" in the title.
$titleLead = $d.Content
$titleLead.Find.Execute("Capstone Project: algo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$bmRange = $d.Range($titleLead.End, $titleLead.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
